# Applies the "Add files via upload" change to 1_3_world_pop.xlsx:
#   - Fills in a new question ("2_" sheet, the 3rd tab) with answer options
#     and comments about average growth of world population 1970-2016.
#   - Makes the "2_" sheet the active/selected tab (was previously "1_").
#   - Updates the answer-range selection on the "4_MultC" sheet.

$wb = $excel.ActiveWorkbook

# The "2_" sheet is the 3rd worksheet in the workbook.
$ws2_ = $wb.Worksheets.Item(3)

# Populate the answer options and their comments first, then the correct/
# label columns, and finally the question text itself -- this mirrors the
# order the strings were authored in and keeps the shared-string table in
# the same sequence as the original edit.
$ws2_.Range("A2").Value = "Higher"
$ws2_.Range("A3").Value = "About the same"
$ws2_.Range("A4").Value = "Lower"
$ws2_.Range("C2").Value = "Yep!  The census and UN curves have a higher slope than our fitted curve during this period."
$ws2_.Range("C3").Value = "Look at the slope of the 3 curves.  How is the slope related to average growth?"
$ws2_.Range("C4").Value = "Look at the slope of the 3 curves.  How is the slope related to average growth?"
$ws2_.Range("A1").Value = "Is the average growth between 1970 and 2016 in the known data higher or lower than the average growth we calculated?"
$ws2_.Range("B1").Value = "Correct"
$ws2_.Range("C1").Value = "Comment"
$ws2_.Range("B2").Value = "Y"
$ws2_.Range("B3").Value = "N"
$ws2_.Range("B4").Value = "N"

# Row heights: header/question row is taller, the 3 answer rows match.
$ws2_.Rows.Item(1).RowHeight = 75
$ws2_.Rows.Item(2).RowHeight = 45
$ws2_.Rows.Item(3).RowHeight = 45
$ws2_.Rows.Item(4).RowHeight = 45

# Make "2_" the selected/active sheet, with C8 as the selected cell.
$ws2_.Activate() | Out-Null
$ws2_.Range("C8").Select() | Out-Null

# Update the stored selection on the "4_MultC" sheet (5th worksheet) to
# match the smaller answer range A1:C6 (previously A1:C10).
$ws4MultC = $wb.Worksheets.Item(5)
$ws4MultC.Range("A1:C6").Select() | Out-Null

# Re-activate "2_" so it remains the active tab when the workbook is saved.
$ws2_.Activate() | Out-Null
